$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''64.455.83'
$ws.Range("E2").Value = '  -5.33%  '
$ws.Range("D3").Value = '''3.342.37'
$ws.Range("E3").Value = '  -7.48%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '''184.04'
$ws.Range("E5").Value = '  -8.95%  '
$ws.Range("D6").Value = '''525.43'
$ws.Range("E6").Value = '  -8.55%  '
$ws.Range("D7").Value = '''0.597'
$ws.Range("E7").Value = '  -3.31%  '
$ws.Range("D8").Value = '''3.337.18'
$ws.Range("E8").Value = '  -7.51%  '
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("D10").Value = '''0.620'
$ws.Range("E10").Value = '  -9.20%  '
$ws.Range("D11").Value = '''58.80'
$ws.Range("E11").Value = '  -2.84%  '
$ws.Range("E12").Value = '  -11.00%  '
$ws.Range("D13").Value = '''0.0000253'
$ws.Range("E13").Value = '  -10.60%  '
$ws.Range("D14").Value = '''9.11'
$ws.Range("E14").Value = '  -12.48%  '
$ws.Range("D15").Value = '''3.892.42'
$ws.Range("E15").Value = '  -7.41%  '
$ws.Range("E16").Value = '  -4.48%  '
$ws.Range("D17").Value = '''3.357.39'
$ws.Range("E17").Value = '  -7.64%  '
$ws.Range("D18").Value = '''64.262.71'
$ws.Range("E18").Value = '  -5.51%  '
$ws.Range("D19").Value = '''17.24'
$ws.Range("E19").Value = '  -10.63%  '
$ws.Range("D20").Value = '''10.92'
$ws.Range("E20").Value = '  -11.27%  '
$ws.Range("D21").Value = '''0.955'
$ws.Range("E21").Value = '  -10.56%  '
$ws.Range("D22").Value = '''372.26'
$ws.Range("E22").Value = '  -8.22%  '
$ws.Range("D23").Value = '''3.70'
$ws.Range("E23").Value = '  -12.07%  '
$ws.Range("D24").Value = '''80.55'
$ws.Range("E24").Value = '  -5.54%  '
$ws.Range("D25").Value = '''10.81'
$ws.Range("E25").Value = '  -16.64%  '
$ws.Range("D26").Value = '''3.81'
$ws.Range("E26").Value = '  -3.80%  '
$ws.Range("D27").Value = '''5.98'
$ws.Range("E27").Value = '  -2.25%  '
$ws.Range("E28").Value = '  -9.11%  '
$ws.Range("D29").Value = '''11.26'
$ws.Range("E29").Value = '  -10.46%  '
$ws.Range("E30").Value = '  -10.04%  '
$ws.Range("D31").Value = '''28.74'
$ws.Range("E31").Value = '  -9.00%  '
$ws.Range("D32").Value = '''652.97'
$ws.Range("E32").Value = '  -3.84%  '
$ws.Range("E33").Value = '  -12.72%  '
$ws.Range("E34").Value = '  -9.01%  '
$ws.Range("D35").Value = '''59.79'
$ws.Range("E35").Value = '  -6.07%  '
$ws.Range("E36").Value = '  -8.82%  '
$ws.Range("E37").Value = '  -0.04%  '
$ws.Range("D38").Value = '''36.22'
$ws.Range("E38").Value = '  -13.13%  '
$ws.Range("E39").Value = '  -8.37%  '
$ws.Range("D40").Value = '''1.00'
$ws.Range("E40").Value = '  +0.12%  '
$ws.Range("E41").Value = '  -7.75%  '
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '''2.807.02'
$ws.Range("E42").Value = '  -12.15%  '
$ws.Range("B43").Value = 'ThetaToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D43").Value = '''2.72'
$ws.Range("E43").Value = '  -14.72%  '
$ws.Range("E44").Value = '  -18.11%  '
$ws.Range("B45").Value = 'WEMIXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").Value = '''2.61'
$ws.Range("E45").Value = '  -7.89%  '
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").Value = '''0.0389'
$ws.Range("E46").Value = '  -5.89%  '
$ws.Range("D47").Value = '''2.33'
$ws.Range("E47").Value = '  -13.92%  '
$ws.Range("E48").Value = '  -5.46%  '
$ws.Range("D49").Value = '''135.13'
$ws.Range("E49").Value = '  -3.12%  '
$ws.Range("D50").Value = '''2.65'
$ws.Range("E50").Value = '  -2.66%  '
$ws.Range("D51").Value = '''2.31'
$ws.Range("E51").Value = '  -18.58%  '
